{"js": "// Replace the 25 \"AxB=C\" multiplication-fact cells with their new values,\n// in document order. Each old value is unique in the document, so a\n// matchCase whole-string search-and-replace is unambiguous.\nconst replacements = [\n  [\"78\u00d747=3666\", \"97\u00d743=4171\"],\n  [\"44\u00d729=1276\", \"66\u00d748=3168\"],\n  [\"77\u00d767=5159\", \"69\u00d796=6624\"],\n  [\"17\u00d711=187\", \"54\u00d733=1782\"],\n  [\"14\u00d790=1260\", \"86\u00d758=4988\"],\n  [\"37\u00d727=999\", \"68\u00d789=6052\"],\n  [\"97\u00d742=4074\", \"33\u00d772=2376\"],\n  [\"98\u00d769=6762\", \"67\u00d738=2546\"],\n  [\"41\u00d786=3526\", \"78\u00d745=3510\"],\n  [\"58\u00d713=754\", \"69\u00d717=1173\"],\n  [\"28\u00d785=2380\", \"72\u00d782=5904\"],\n  [\"66\u00d792=6072\", \"17\u00d787=1479\"],\n  [\"32\u00d784=2688\", \"94\u00d768=6392\"],\n  [\"47\u00d795=4465\", \"83\u00d744=3652\"],\n  [\"77\u00d796=7392\", \"57\u00d716=912\"],\n  [\"86\u00d767=5762\", \"94\u00d743=4042\"],\n  [\"94\u00d724=2256\", \"59\u00d739=2301\"],\n  [\"47\u00d725=1175\", \"84\u00d761=5124\"],\n  [\"38\u00d781=3078\", \"83\u00d783=6889\"],\n  [\"32\u00d759=1888\", \"69\u00d763=4347\"],\n  [\"75\u00d747=3525\", \"14\u00d715=210\"],\n  [\"80\u00d758=4640\", \"75\u00d744=3300\"],\n  [\"61\u00d788=5368\", \"62\u00d783=5146\"],\n  [\"39\u00d778=3042\", \"34\u00d779=2686\"],\n  [\"56\u00d768=3808\", \"39\u00d783=3237\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"AxB=C\" multiplication-fact cells with their new values,\n# in document order. Each old value is unique in the document, so a\n# MatchCase whole-document Find/Replace is unambiguous and safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"78\u00d747=3666\", \"97\u00d743=4171\"),\n    @(\"44\u00d729=1276\", \"66\u00d748=3168\"),\n    @(\"77\u00d767=5159\", \"69\u00d796=6624\"),\n    @(\"17\u00d711=187\", \"54\u00d733=1782\"),\n    @(\"14\u00d790=1260\", \"86\u00d758=4988\"),\n    @(\"37\u00d727=999\", \"68\u00d789=6052\"),\n    @(\"97\u00d742=4074\", \"33\u00d772=2376\"),\n    @(\"98\u00d769=6762\", \"67\u00d738=2546\"),\n    @(\"41\u00d786=3526\", \"78\u00d745=3510\"),\n    @(\"58\u00d713=754\", \"69\u00d717=1173\"),\n    @(\"28\u00d785=2380\", \"72\u00d782=5904\"),\n    @(\"66\u00d792=6072\", \"17\u00d787=1479\"),\n    @(\"32\u00d784=2688\", \"94\u00d768=6392\"),\n    @(\"47\u00d795=4465\", \"83\u00d744=3652\"),\n    @(\"77\u00d796=7392\", \"57\u00d716=912\"),\n    @(\"86\u00d767=5762\", \"94\u00d743=4042\"),\n    @(\"94\u00d724=2256\", \"59\u00d739=2301\"),\n    @(\"47\u00d725=1175\", \"84\u00d761=5124\"),\n    @(\"38\u00d781=3078\", \"83\u00d783=6889\"),\n    @(\"32\u00d759=1888\", \"69\u00d763=4347\"),\n    @(\"75\u00d747=3525\", \"14\u00d715=210\"),\n    @(\"80\u00d758=4640\", \"75\u00d744=3300\"),\n    @(\"61\u00d788=5368\", \"62\u00d783=5146\"),\n    @(\"39\u00d778=3042\", \"34\u00d779=2686\"),\n    @(\"56\u00d768=3808\", \"39\u00d783=3237\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\nWrite-Output \"done\"\n"}
